$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final player data for rows 2-18 (after removing Jusuf Nurkic and reordering)
$data = New-Object 'object[,]' 17,3
$data[0,0] = 'Tyler Herro'
$data[0,1] = 'PG,SG'
$data[0,2] = 'Miami Heat'
$data[1,0] = 'Luka Doncic'
$data[1,1] = 'PG,SG'
$data[1,2] = 'Dallas Mavericks'
$data[2,0] = 'Bilal Coulibaly'
$data[2,1] = 'SG,SF'
$data[2,2] = 'Washington Wizards'
$data[3,0] = 'Bennedict Mathurin'
$data[3,1] = 'SG,SF'
$data[3,2] = 'Indiana Pacers'
$data[4,0] = 'Michael Porter Jr.'
$data[4,1] = 'SF,PF'
$data[4,2] = 'Denver Nuggets'
$data[5,0] = 'Miles Bridges'
$data[5,1] = 'SF,PF'
$data[5,2] = 'Charlotte Hornets'
$data[6,0] = 'Mikal Bridges'
$data[6,1] = 'SG,SF,PF'
$data[6,2] = 'New York Knicks'
$data[7,0] = 'Amen Thompson'
$data[7,1] = 'SG,SF'
$data[7,2] = 'Houston Rockets'
$data[8,0] = 'Evan Mobley'
$data[8,1] = 'PF,C'
$data[8,2] = 'Cleveland Cavaliers'
$data[9,0] = 'Nikola Vucevic'
$data[9,1] = 'PF,C'
$data[9,2] = 'Chicago Bulls'
$data[10,0] = 'Josh Giddey'
$data[10,1] = 'PG,SG,SF'
$data[10,2] = 'Chicago Bulls'
$data[11,0] = 'Brook Lopez'
$data[11,1] = 'C'
$data[11,2] = 'Milwaukee Bucks'
$data[12,0] = 'De''Aaron Fox'
$data[12,1] = 'PG'
$data[12,2] = 'Sacramento Kings'
$data[13,0] = 'Buddy Hield'
$data[13,1] = 'SG,SF'
$data[13,2] = 'Golden State Warriors'
$data[14,0] = 'Deandre Ayton'
$data[14,1] = 'C'
$data[14,2] = 'Portland Trail Blazers'
$data[15,0] = 'Scottie Barnes'
$data[15,1] = 'SG,SF,PF'
$data[15,2] = 'Toronto Raptors'
$data[16,0] = 'DeMar DeRozan'
$data[16,1] = 'SF,PF'
$data[16,2] = 'Sacramento Kings'

$ws.Range("A2:C18").Value = $data

# Remove the now-extra row 19
$ws.Rows.Item(19).Delete()

$ws.Range("A1").Select()
